# Apply the "alpha_non_zero" regeneration of this non-convex experiment
# workbook: new numeric values for the leader/follower restriction
# expressions, the modified point, and the bf/BF/alpha vectors.
#
# Many of these cells hold numbers that were authored as *text* (the
# original workbook stores them as shared strings, not as <v> numeric
# cells). Plain `Range.Value = "-3.89"` would let Excel's type-sniffer
# convert that back into a real number, so for every such cell we
# temporarily force Text format, assign the literal string, then reset
# the style back to Normal (this clears the format back to the original
# General/style-0 look without leaving a stray NumberFormat behind).
#
# NOTE: worksheets are addressed by their 1-based tab index rather than
# by name — this workbook has both "Vector_bf" and "Vector_BF" tabs, and
# name-based `Worksheets.Item(...)` lookup is case-insensitive, so it
# would resolve both names to the same sheet.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, [string]$addr, [string]$val)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# --- 2: Restricciones_del_lider -------------------------------------------
$ws = $wb.Worksheets.Item(2)

Set-TextValue $ws "A2" "1.9399999999999995 - x"
Set-TextValue $ws "B2" "-2.9399999999999995"
Set-TextValue $ws "D2" "0.32"

Set-TextValue $ws "A3" "-1.9399999999999997 + x"
Set-TextValue $ws "B3" "0.9399999999999997"
Set-TextValue $ws "D3" "0.02"

Set-TextValue $ws "A4" "35.63239999999999 + x - y - 9(x^2)"
Set-TextValue $ws "B4" "-34.63239999999999"
Set-TextValue $ws "D4" "0.44"

# --- 3: Restricciones_del_follower -----------------------------------------
$ws = $wb.Worksheets.Item(3)

Set-TextValue $ws "A2" "17.688724378109445 - 10.108736318407956y + (-0.5 + x)*(y^2)"
Set-TextValue $ws "B2" "-17.688724378109445"
Set-TextValue $ws "D2" "0.36"
Set-TextValue $ws "E2" "0"
Set-TextValue $ws "F2" "1.1"

Set-TextValue $ws "A3" "0"
Set-TextValue $ws "B3" "-1"
Set-TextValue $ws "D3" "0.43"
Set-TextValue $ws "E3" "0"
Set-TextValue $ws "F3" "8.7"

Set-TextValue $ws "A4" "-16.260845771144275 + 4.278606965174129y"
Set-TextValue $ws "B4" "14.830845771144274"
Set-TextValue $ws "D4" "0.43"
Set-TextValue $ws "E4" "0"
Set-TextValue $ws "F4" "8.6"

# --- 4: Punto_modificado -----------------------------------------------------
$ws = $wb.Worksheets.Item(4)

Set-TextValue $ws "A2" "1.9399999999999997"
Set-TextValue $ws "B2" "3.6999999999999993"

# --- 5: Vector_bf ----------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)

Set-TextValue $ws "A2" "-3.0368159203980096"

# --- 6: Vector_BF ------------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)

Set-TextValue $ws "A2" "5.99016"
Set-TextValue $ws "A3" "-1.6439679999999997"

# --- 7: Vector_Alpha (genuine numeric cell, unlike the ones above) -----------------
$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = 2.0100000000000002
